# Auto-generated edit script: updates Kraken_Profits leve-profit calculations
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 5602
$ws.Range("I4").Value = 1377.5
$ws.Range("J4").Value = 22500
$ws.Range("K4").Value = 1377.5
$ws.Range("L4").Value = 22500
$ws.Range("M4").Value = -1263.5

$ws.Range("H5").Value = 136.42857
$ws.Range("I5").Value = 139.23077
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 139.23077
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -24.23077000000001
$ws.Range("N5").Value = -330

$ws.Range("H19").Value = 406.85715
$ws.Range("I19").Value = 300
$ws.Range("J19").Value = 449.6
$ws.Range("K19").Value = 300
$ws.Range("L19").Value = 449.6
$ws.Range("M19").Value = -125
$ws.Range("N19").Value = -799.6

$ws.Range("H40").Value = 8019.3477
$ws.Range("I40").Value = 4099.4
$ws.Range("J40").Value = 9108.223
$ws.Range("K40").Value = 4099.4
$ws.Range("L40").Value = 9108.223
$ws.Range("M40").Value = -3924.4
$ws.Range("N40").Value = -9458.223

$ws.Range("H43").Value = 6999.5
$ws.Range("I43").Value = 6999.5
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 6999.5
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -6930.5

$ws.Range("H80").Value = 5231.727
$ws.Range("I80").Value = 600
$ws.Range("J80").Value = 9091.5
$ws.Range("K80").Value = 1800
$ws.Range("L80").Value = 27274.5
$ws.Range("M80").Value = -802

$ws.Range("H83").Value = 5231.727
$ws.Range("I83").Value = 600
$ws.Range("J83").Value = 9091.5
$ws.Range("K83").Value = 5400
$ws.Range("L83").Value = 81823.5
$ws.Range("M83").Value = -408

$ws.Range("H86").Value = 1848
$ws.Range("I86").Value = 1956
$ws.Range("J86").Value = 1200
$ws.Range("K86").Value = 1956
$ws.Range("L86").Value = 1200
$ws.Range("M86").Value = -833

$ws.Range("H89").Value = 1848
$ws.Range("I89").Value = 1956
$ws.Range("J89").Value = 1200
$ws.Range("K89").Value = 9780
$ws.Range("L89").Value = 6000
$ws.Range("M89").Value = -4164

$ws.Range("H93").Value = 38800
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 38800
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 38800
$ws.Range("N93").Value = -43792

$ws.Range("H100").Value = 1907.4286
$ws.Range("I100").Value = 1907.4286
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1907.4286
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1366.4286

$ws.Range("H130").Value = 75421.71000000001
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 75421.71000000001
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 75421.71000000001
$ws.Range("N130").Value = -85461.71000000001

$ws.Range("H132").Value = 1902.5652
$ws.Range("I132").Value = 1902.5652
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5707.6956
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3177.6956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 591
$ws.Range("I97").Value = 607
$ws.Range("J97").Value = 511
$ws.Range("K97").Value = 607
$ws.Range("L97").Value = 511
$ws.Range("M97").Value = -111

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4667
$ws.Range("I107").Value = 2000.5
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 2000.5
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = -80.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1703.6666
$ws.Range("I16").Value = 1305.5
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 1305.5
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -1018.5

$ws.Range("H31").Value = 2476.4285
$ws.Range("I31").Value = 803.4
$ws.Range("J31").Value = 3997.3635
$ws.Range("K31").Value = 803.4
$ws.Range("L31").Value = 3997.3635
$ws.Range("M31").Value = -508.4
$ws.Range("N31").Value = -4587.363499999999

$ws.Range("H34").Value = 2476.4285
$ws.Range("I34").Value = 803.4
$ws.Range("J34").Value = 3997.3635
$ws.Range("K34").Value = 803.4
$ws.Range("L34").Value = 3997.3635
$ws.Range("M34").Value = -601.4
$ws.Range("N34").Value = -4401.363499999999

$ws.Range("H57").Value = 6000
$ws.Range("I57").Value = 5000
$ws.Range("J57").Value = 6500
$ws.Range("K57").Value = 5000
$ws.Range("L57").Value = 6500
$ws.Range("M57").Value = -4440
$ws.Range("N57").Value = -7620

$ws.Range("H93").Value = 21975
$ws.Range("I93").Value = 21975
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 21975
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -20103

$ws.Range("H113").Value = 1703.6666
$ws.Range("I113").Value = 1305.5
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1305.5
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 864.5

$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 750463.1
$ws.Range("I4").Value = 714815.1
$ws.Range("J4").Value = 999999
$ws.Range("K4").Value = 2144445.3
$ws.Range("L4").Value = 2999997
$ws.Range("M4").Value = -2144333.3

$ws.Range("H11").Value = 7143168
$ws.Range("I11").Value = 12500349
$ws.Range("J11").Value = 260.33334
$ws.Range("K11").Value = 37501047
$ws.Range("L11").Value = 781.0000200000001
$ws.Range("M11").Value = -37500907
$ws.Range("N11").Value = -1061.00002

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 25000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 25000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 25000
$ws.Range("N40").Value = -25302
$ws.Range("M40").ClearContents()

$ws.Range("H92").Value = 8217.1
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 8217.1
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 8217.1
$ws.Range("N92").Value = -11961.1
$ws.Range("M92").ClearContents()

$ws.Range("H93").Value = 65000
$ws.Range("I93").Value = 55000
$ws.Range("J93").Value = 70000
$ws.Range("K93").Value = 55000
$ws.Range("L93").Value = 70000
$ws.Range("M93").Value = -53128
$ws.Range("N93").Value = -73744

$ws.Range("H107").Value = 1162.8334
$ws.Range("I107").Value = 200.66667
$ws.Range("J107").Value = 2125
$ws.Range("K107").Value = 200.66667
$ws.Range("L107").Value = 2125
$ws.Range("M107").Value = 1719.33333
$ws.Range("N107").Value = -5965

$ws.Range("H134").Value = 99999
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 99999
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 299997
$ws.Range("N134").Value = -305067

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 11166.667
$ws.Range("I4").Value = 3500
$ws.Range("J4").Value = 15000
$ws.Range("K4").Value = 3500
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = -3387
$ws.Range("N4").Value = -15226

$ws.Range("H5").Value = 9666.666999999999
$ws.Range("I5").Value = 12000
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = -11887
$ws.Range("N5").Value = -5226

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()

$ws.Range("H28").Value = 11166.667
$ws.Range("I28").Value = 3500
$ws.Range("J28").Value = 15000
$ws.Range("K28").Value = 3500
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = -3268
$ws.Range("N28").Value = -15464

$ws.Range("H31").Value = 6267.1113
$ws.Range("I31").Value = 3767.1428
$ws.Range("J31").Value = 15017
$ws.Range("K31").Value = 3767.1428
$ws.Range("L31").Value = 15017
$ws.Range("M31").Value = -3519.1428
$ws.Range("N31").Value = -15513

$ws.Range("H37").Value = 11166.667
$ws.Range("I37").Value = 3500
$ws.Range("J37").Value = 15000
$ws.Range("K37").Value = 3500
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = -3393
$ws.Range("N37").Value = -15214

$ws.Range("H68").Value = 5416.5
$ws.Range("I68").Value = 2499.8
$ws.Range("J68").Value = 20000
$ws.Range("K68").Value = 2499.8
$ws.Range("L68").Value = 20000
$ws.Range("M68").Value = -1750.8

$ws.Range("H71").Value = 5416.5
$ws.Range("I71").Value = 2499.8
$ws.Range("J71").Value = 20000
$ws.Range("K71").Value = 12499
$ws.Range("L71").Value = 100000
$ws.Range("M71").Value = -8755

$ws.Range("H82").Value = 2165.3076
$ws.Range("I82").Value = 1683.3334
$ws.Range("J82").Value = 3249.75
$ws.Range("K82").Value = 1683.3334
$ws.Range("L82").Value = 3249.75
$ws.Range("M82").Value = -1322.3334

$ws.Range("H85").Value = 2165.3076
$ws.Range("I85").Value = 1683.3334
$ws.Range("J85").Value = 3249.75
$ws.Range("K85").Value = 1683.3334
$ws.Range("L85").Value = 3249.75
$ws.Range("M85").Value = -435.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 30025
$ws.Range("I40").Value = 30025
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 30025
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -29876

$ws.Range("H57").Value = 166855.67
$ws.Range("I57").Value = 100177
$ws.Range("J57").Value = 200195
$ws.Range("K57").Value = 100177
$ws.Range("L57").Value = 200195
$ws.Range("M57").Value = -99423
$ws.Range("N57").Value = -201703

$ws.Range("H92").Value = 40000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 40000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992

$ws.Range("H107").Value = 3400.4
$ws.Range("I107").Value = 1667.3334
$ws.Range("J107").Value = 6000
$ws.Range("K107").Value = 5002.0002
$ws.Range("L107").Value = 18000
$ws.Range("M107").Value = -3082.0002
$ws.Range("N107").Value = -21840

$ws.Range("H113").Value = 708.26666
$ws.Range("I113").Value = 727.0833
$ws.Range("J113").Value = 633
$ws.Range("K113").Value = 2181.2499
$ws.Range("L113").Value = 1899
$ws.Range("M113").Value = -11.2498999999998

$ws.Range("H122").Value = 2936.353
$ws.Range("I122").Value = 3355.2307
$ws.Range("J122").Value = 1575
$ws.Range("K122").Value = 10065.6921
$ws.Range("L122").Value = 4725
$ws.Range("M122").Value = -7615.6921

$ws.Range("H136").Value = 3978.4285
$ws.Range("I136").Value = 3371.1
$ws.Range("J136").Value = 5496.75
$ws.Range("K136").Value = 10113.3
$ws.Range("L136").Value = 16490.25
$ws.Range("M136").Value = -7563.299999999999

Write-Output "Applied Kraken_Profits updates"
